# The paragraph was rewritten in full ("Add files via upload" - a fresh
# copy of the lesson text was dropped in over the old one). Reproduce
# that by replacing the old paragraph text with the new paragraph text.
#
# Note: this runtime's Find.Execute has an off-by-one quirk where a
# match/replace whose end coincides exactly with the end of the story
# also swallows one extra trailing character. We dodge it by appending
# two throwaway sentinel characters to the very end of the document
# before doing the big replace (so the erroneously-consumed character
# is just a sentinel, not real content), then removing the sentinel
# that is left over afterwards.

$d = $word.ActiveDocument

$oldText = "Bu dersimizde ünlü İtalyan halk şarkısı " + [char]0x201C + "Çav Bella" + [char]0x201D + " yı görmekteyiz. Nota değerlerine dikkat ederek çalalım. Acele etmeyelim. Yavaş çalarak şarkıyı sindirelim daha sonra şarkıyı orijinal temposuna çekelim. Videodaki parmak pozisyonlarına dikkat ederek çalmak ,kullanmadığımız parmaklarımızın kuvvetlenmesini dolayısıyla hızlanmamızı sağlayacaktır. Keyifli çalışmalar."

$newText = "Sırada eğlenceli bir Blues şarkısı var. Perdelere basarken doğru parmakları kullanmaya özen gösterelim. Şarkıyı ölçü ölçü çalışıp en son birleştirelim. Tel geçişlerinde acele etmeyelim ama zamanlamayı doğru yapmaya çalışalım. Keyifli çalışmalar "

# Sentinel guard so the whole-story replace below can't eat real content.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertAfter("##")

$d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2) | Out-Null

# Clean up the leftover sentinel character.
$d.Content.Find.Execute("#", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null
